# Course TimeSheet update: new class format/guest details for upcoming sessions.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell content updates ("توضیحات" column D, and a few "موضوع"/"تاریخ" tweaks) ---
$ws.Range("D9").Value2  = "هیبریدی برگزار میشود"
$ws.Range("D10").Value2 = "مجازی برگزار میشود"
$ws.Range("D11").Value2 = "مجازی برگزار میشود"
$ws.Range("D12").Value2 = "قراره که درباره این تاریخ هم با دانشکده مشورت و صحبت کنم"

$ws.Range("E26").Value2 = "مهمان (سینا)"
$ws.Range("E27").Value2 = "مهمان (مسعود)"
$ws.Range("E28").Value2 = "!گردش فناورانه"

$ws.Range("C29").Value2 = "17/3"
$ws.Range("D29").Value2 = "کلاس جبرانی"

# --- Column D needs to be wider now that it holds longer notes ---
$ws.Columns.Item(4).ColumnWidth = 41

# --- Update the view/selection to where the edits were made ---
$ws.Range("D10").Select()
$excel.ActiveWindow.ScrollRow = 12
$excel.ActiveWindow.ScrollColumn = 1

Write-Host "Applied Course TimeSheet updates"
